$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Stimulus, Brexit, Storm:Your Thursday Evening Briefing", "https://www.nytimes.com/2020/12/24/briefing/stimulus-brexit-storm.html"),
    @("UKs Gove says: theres a chance of a no-deal Brexit - Reuters India", "https://www.reuters.com/article/uk-britain-eu-gove-idINKBN28B487"),
    @("UKs Gove says: theres a chance of a no-deal Brexit - Reuters UK", "https://uk.reuters.com/article/uk-britain-eu-gove-idUKKBN28B487"),
    @("EU sets out plans in case Brexit trade talks fail", "https://www.bbc.co.uk/news/world-europe-55259144"),
    @("UKs Gove sees scope for compromise on fish in Brexit talks - Reuters", "https://www.reuters.com/article/us-britain-eu-fish-idUSKBN28J0WU"),
    @("UKs Gove sees scope for compromise on fish in Brexit talks - Reuters", "https://www.reuters.com/article/uk-britain-eu-fish-idUSKBN28J0X0"),
    @("Sterling holds above `$1.33 on Brexit trade deal hopes - Reuters UK", "https://uk.reuters.com/article/uk-britain-sterling-idUKKBN28A10E"),
    @("Stena Line says role in UK delivery of Pfizer vaccine going well - Reuters UK", "https://uk.reuters.com/article/uk-health-coronavirus-britain-stena-idUKKBN28E191"),
    @("UK cabinet to back Johnson over no-deal Brexit - The Times - Reuters", "https://www.reuters.com/article/uk-britain-eu-trade-idUSKBN28G010"),
    @("FACTBOX-Higher prices, traffic jams: likely impact of Brexit impasse on businesses - Reuters", "https://www.reuters.com/article/britain-eu-companies-idUSL8N2IN27C"),
    @("Sterling stabilises after Brexit deal hopes rebound - Reuters", "https://www.reuters.com/article/britain-sterling-idUSKBN28P0XT"),
    @('EUs Barnier says "crucial" moment in UK trade talks - Reuters', "https://www.reuters.com/article/britain-eu-barnier-idUSKBN28U0J8"),
    @("Factbox-Higher prices, traffic jams: likely impact of Brexit impasse on businesses - Reuters UK", "https://uk.reuters.com/article/uk-britain-eu-companies-idUKKBN28K1P2"),
    @("EUs Barnier tells member states latest UK fish offer unacceptable - Reuters Canada", "https://www.reuters.com/article/us-britain-eu-barnier-fish-idUSKBN28W214"),
    @("Rifts over fishing rights still blocking EU-UK trade deal - Reuters India", "https://in.reuters.com/article/us-britain-eu-idINKBN28W13M"),
    @("Rifts over fishing rights still blocking EU-UK trade deal - Reuters", "https://www.reuters.com/article/uk-britain-eu-barnier-idUSKBN28W13U"),
    @("PRECIOUS-Gold gains on subdued dollar after Brexit trade deal - Reuters UK", "https://in.reuters.com/article/global-precious-idUKL4N2J4274"),
    @("Smoother glide path to Brexit trade deal after Northern Ireland accord, UKs Gove says - Reuters", "https://www.reuters.com/article/britain-eu-gove-idUSKBN28J0PS"),
    @("PRECIOUS-Gold gains on subdued dollar after Brexit trade deal - Reuters UK", "https://uk.reuters.com/article/global-precious-idUKL4N2J4274"),
    @("Stena Line says role in UK delivery of Pfizer vaccine going well - Reuters", "https://www.reuters.com/article/uk-health-coronavirus-britain-stena-idUSKBN28E191")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
